$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.892.88"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.809.12"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4968"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3927"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +23.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.097"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.418"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "1.811.63"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.263"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001125"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06645"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.903"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "27.950.18"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.248"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("D27").Value = "2.020.00"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.378"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1061"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.030"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.542"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.614"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06703"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.888"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2125"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6147"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.147"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.294"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5864"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.693"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.924"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06752"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.16%  "
